# Updating filtered feeds from workflow
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filtered Feeds")

# New feed rows appended to the bottom of the "Filtered Feeds" table
# (rows 16-19), matching the A:link / B:keyword / C:title layout used
# by the existing rows.
$links = @(
    "https://www.sciencedaily.com/releases/2024/12/241210163409.htm",
    "https://www.sciencedaily.com/releases/2024/01/240114202019.htm",
    "https://www.sciencedaily.com/releases/2019/10/191028104212.htm",
    "https://www.sciencedaily.com/releases/2018/04/180426141507.htm"
)
$keywords = @(
    "smoldering multiple myeloma",
    "MGUS",
    "smoldering multiple myeloma",
    "MGUS"
)
$titles = @(
    "Research shows new treatment could delay cancer progression in patients with high-risk smoldering multiple myeloma",
    "Obesity linked to detection of blood cancer precursor",
    "Lenalidomide may delay onset of myeloma-related bone, organ damage",
    "Blood cancer precursor found in 9/11 firefighters"
)

$startRow = 16
$sampleStyle = $ws.Cells.Item($startRow - 1, 1).Style

for ($i = 0; $i -lt $links.Length; $i++) {
    $row = $startRow + $i
    $linkCell = $ws.Cells.Item($row, 1)

    $linkCell.Value = $links[$i]
    $ws.Hyperlinks.Add($linkCell, $links[$i])
    $linkCell.Style = $sampleStyle

    $ws.Cells.Item($row, 2).Value = $keywords[$i]
    $ws.Cells.Item($row, 3).Value = $titles[$i]
}
